$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin name / link swaps (rows 6-7 and 35-36)
$ws.Range("B6").Value = 'USDC'
$ws.Range("C6").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("B36").Value = 'Maker'
$ws.Range("C36").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'

# Price / Volume(1h) updates - force text format to match source data (values like "1.01", "26.424.48" must stay text, not become numbers/dates)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.424.48'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.05%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.624.10'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.51%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.58%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.69'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.01'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.55%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.494'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.84%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.249'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.85%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0619'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.13%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.91'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.74%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0839'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.70%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.853.72'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.22%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.668.60'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.30%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.10'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.31%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.519'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.99%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.77'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.11%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.449.77'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.83%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0₃0737'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.44%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '214.88'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.61%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.48%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.29'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.97%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.18'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.58%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.28'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.22%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +4.55%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '148.22'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.69%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.43%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.119'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.82'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.46%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.53'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.07%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.36%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.68%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.30'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.20%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.92'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.36%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.49'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.36%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.38'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.45%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.214.73'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +4.20%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +4.07%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.54%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.792'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.76%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.499'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.72%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.80%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.793'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.83%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.99%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.763.75'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.17%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.91'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.62%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.56'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.64%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '54.56'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.21%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0₆0102'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.81%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.15%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.58'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.19%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.409'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.37%  '
